# Apply "GOOD DATA FOR 8 GANTT CHARTS" edit:
#  - Rename header G1 from END_PM to "End PM"
#  - Convert G2:G27 numeric End PM values into formatted inline-string
#    text cells ("End PM:   58.583", matching the existing Beg PM format)
#  - Update the SQL text on the "SQL" worksheet so both SELECT branches
#    wrap end_pm with the same to_char(...) formatting as beg_pm

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Export Worksheet" ---
$ws = $wb.Worksheets.Item("Export Worksheet")

# Header rename
$ws.Range("G1").Value = "End PM"

# Data rows 2..27: replace numeric End PM with formatted text, matching
# the "Beg PM:" style already used in column F (label + value right
# justified in an 8-character field with 3 decimals).
for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 7)
    $numValue = $cell.Value2
    $formatted = "End PM: {0,8:N3}" -f $numValue
    $cell.Value = $formatted
    $cell.ClearFormats()
}

# --- Sheet 2: "SQL" ---
$sqlWs = $wb.Worksheets.Item("SQL")

$newSql = "select a.ea, a.treatment, a.county, a.route, a.year, ('Beg PM: ' || to_char(a.beg_pm, 990.999)) as ""Beg PM"", ('End PM: ' || to_char(a.end_pm, 990.999)) as ""End PM"", (a.end_pm-a.beg_pm) as length, a.budget_group from s1383currentr a `nwhere a.county = 'SON'  `nunion  `nselect b.ea, b.treatment, b.county, b.route, b.year, ('Beg PM: ' || to_char(b.beg_pm, 990.999)) as ""Beg PM"",  ('End PM: ' || to_char(b.end_pm, 990.999)) as ""End PM"",  (b.end_pm-b.beg_pm) as length, b.budget_group from s1383historyr b `nwhere b.county = 'SON'  `norder by year"

$sqlWs.Range("A2").Value = $newSql
